$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("checkaxiom")
$ws2 = $wb.Worksheets.Item("aei")

# --- Update existing shared text labels (casing fix) on header rows ---
$ws1.Range("A1").Value = "eWGARP"
$ws1.Range("B1").Value = "eWGARP"
$ws1.Range("C1").Value = "eWGARP"
$ws1.Range("D1").Value = "eWARP"
$ws1.Range("E1").Value = "eWARP"
$ws1.Range("F1").Value = "eWARP"
$ws1.Range("G1").Value = "eGARP"
$ws1.Range("H1").Value = "eGARP"
$ws1.Range("I1").Value = "eGARP"
$ws1.Range("J1").Value = "eSARP"
$ws1.Range("K1").Value = "eSARP"
$ws1.Range("L1").Value = "eSARP"
$ws1.Range("M1").Value = "eHARP"
$ws1.Range("N1").Value = "eHARP"
$ws1.Range("O1").Value = "eHARP"
$ws1.Range("P1").Value = "eCM"
$ws1.Range("Q1").Value = "eCM"
$ws1.Range("R1").Value = "eCM"

$ws2.Range("A1").Value = "eWGARP"
$ws2.Range("B1").Value = "eWARP"
$ws2.Range("C1").Value = "eGARP"
$ws2.Range("D1").Value = "eSARP"
$ws2.Range("E1").Value = "eHARP"
$ws2.Range("F1").Value = "eCM"

# --- Special-case updates where an existing metric column value also changed ---
$ws1.Range("K5").Value = 4
$ws1.Range("L5").Value = 7.14
$ws1.Range("H40").Value = 8
$ws1.Range("I40").Value = 14.29
$ws1.Range("K40").Value = 9
$ws1.Range("L40").Value = 16.07
$ws1.Range("H42").Value = 8
$ws1.Range("I42").Value = 14.29
$ws1.Range("K42").Value = 11
$ws1.Range("L42").Value = 19.64
$ws1.Range("K43").Value = 2
$ws1.Range("L43").Value = 3.57
$ws1.Range("K49").Value = 2
$ws1.Range("L49").Value = 3.57
$ws1.Range("H63").Value = 4
$ws1.Range("I63").Value = 7.14
$ws1.Range("K63").Value = 5
$ws1.Range("L63").Value = 8.93
$ws1.Range("K74").Value = 2
$ws1.Range("L74").Value = 3.57
$ws1.Range("K89").Value = 2
$ws1.Range("L89").Value = 3.57
$ws1.Range("H92").Value = 2
$ws1.Range("I92").Value = 3.57
$ws1.Range("K92").Value = 2
$ws1.Range("L92").Value = 3.57
$ws1.Range("K106").Value = 3
$ws1.Range("L106").Value = 5.36
$ws1.Range("K128").Value = 4
$ws1.Range("L128").Value = 7.14
$ws1.Range("K139").Value = 2
$ws1.Range("L139").Value = 3.57
$ws1.Range("K141").Value = 2
$ws1.Range("L141").Value = 3.57
# --- Add new eSGARP column group (S,T,U) to sheet1 ---
$ws1.Range("S1").Value = "eSGARP"
$ws1.Range("T1").Value = "eSGARP"
$ws1.Range("U1").Value = "eSGARP"
$ws1.Range("S2").Value = "Pass"
$ws1.Range("T2").Value = "Violations"
$ws1.Range("U2").Value = "Violations_frac"

$ws1.Range("S3").Value = 1
$ws1.Range("T3").Value = 0
$ws1.Range("U3").Value = 0
$ws1.Range("S4").Value = 0
$ws1.Range("T4").Value = 39
$ws1.Range("U4").Value = 60.94
$ws1.Range("S5").Value = 0
$ws1.Range("T5").Value = 31
$ws1.Range("U5").Value = 48.44
$ws1.Range("S6").Value = 0
$ws1.Range("T6").Value = 25
$ws1.Range("U6").Value = 39.06
$ws1.Range("S7").Value = 0
$ws1.Range("T7").Value = 39
$ws1.Range("U7").Value = 60.94
$ws1.Range("S8").Value = 1
$ws1.Range("T8").Value = 0
$ws1.Range("U8").Value = 0
$ws1.Range("S9").Value = 1
$ws1.Range("T9").Value = 0
$ws1.Range("U9").Value = 0
$ws1.Range("S10").Value = 0
$ws1.Range("T10").Value = 39
$ws1.Range("U10").Value = 60.94
$ws1.Range("S11").Value = 1
$ws1.Range("T11").Value = 0
$ws1.Range("U11").Value = 0
$ws1.Range("S12").Value = 1
$ws1.Range("T12").Value = 0
$ws1.Range("U12").Value = 0
$ws1.Range("S13").Value = 1
$ws1.Range("T13").Value = 0
$ws1.Range("U13").Value = 0
$ws1.Range("S14").Value = 1
$ws1.Range("T14").Value = 0
$ws1.Range("U14").Value = 0
$ws1.Range("S15").Value = 1
$ws1.Range("T15").Value = 0
$ws1.Range("U15").Value = 0
$ws1.Range("S16").Value = 0
$ws1.Range("T16").Value = 39
$ws1.Range("U16").Value = 60.94
$ws1.Range("S17").Value = 1
$ws1.Range("T17").Value = 0
$ws1.Range("U17").Value = 0
$ws1.Range("S18").Value = 0
$ws1.Range("T18").Value = 8
$ws1.Range("U18").Value = 12.5
$ws1.Range("S19").Value = 1
$ws1.Range("T19").Value = 0
$ws1.Range("U19").Value = 0
$ws1.Range("S20").Value = 1
$ws1.Range("T20").Value = 0
$ws1.Range("U20").Value = 0
$ws1.Range("S21").Value = 0
$ws1.Range("T21").Value = 39
$ws1.Range("U21").Value = 60.94
$ws1.Range("S22").Value = 0
$ws1.Range("T22").Value = 5
$ws1.Range("U22").Value = 7.81
$ws1.Range("S23").Value = 0
$ws1.Range("T23").Value = 39
$ws1.Range("U23").Value = 60.94
$ws1.Range("S24").Value = 0
$ws1.Range("T24").Value = 37
$ws1.Range("U24").Value = 57.81
$ws1.Range("S25").Value = 1
$ws1.Range("T25").Value = 0
$ws1.Range("U25").Value = 0
$ws1.Range("S26").Value = 0
$ws1.Range("T26").Value = 25
$ws1.Range("U26").Value = 39.06
$ws1.Range("S27").Value = 1
$ws1.Range("T27").Value = 0
$ws1.Range("U27").Value = 0
$ws1.Range("S28").Value = 0
$ws1.Range("T28").Value = 2
$ws1.Range("U28").Value = 3.13
$ws1.Range("S29").Value = 1
$ws1.Range("T29").Value = 0
$ws1.Range("U29").Value = 0
$ws1.Range("S30").Value = 1
$ws1.Range("T30").Value = 0
$ws1.Range("U30").Value = 0
$ws1.Range("S31").Value = 0
$ws1.Range("T31").Value = 25
$ws1.Range("U31").Value = 39.06
$ws1.Range("S32").Value = 0
$ws1.Range("T32").Value = 14
$ws1.Range("U32").Value = 21.88
$ws1.Range("S33").Value = 0
$ws1.Range("T33").Value = 2
$ws1.Range("U33").Value = 3.13
$ws1.Range("S34").Value = 0
$ws1.Range("T34").Value = 39
$ws1.Range("U34").Value = 60.94
$ws1.Range("S35").Value = 0
$ws1.Range("T35").Value = 2
$ws1.Range("U35").Value = 3.13
$ws1.Range("S36").Value = 0
$ws1.Range("T36").Value = 39
$ws1.Range("U36").Value = 60.94
$ws1.Range("S37").Value = 0
$ws1.Range("T37").Value = 39
$ws1.Range("U37").Value = 60.94
$ws1.Range("S38").Value = 1
$ws1.Range("T38").Value = 0
$ws1.Range("U38").Value = 0
$ws1.Range("S39").Value = 0
$ws1.Range("T39").Value = 6
$ws1.Range("U39").Value = 9.38
$ws1.Range("S40").Value = 0
$ws1.Range("T40").Value = 10
$ws1.Range("U40").Value = 15.63
$ws1.Range("S41").Value = 1
$ws1.Range("T41").Value = 0
$ws1.Range("U41").Value = 0
$ws1.Range("S42").Value = 0
$ws1.Range("T42").Value = 29
$ws1.Range("U42").Value = 45.31
$ws1.Range("S43").Value = 0
$ws1.Range("T43").Value = 4
$ws1.Range("U43").Value = 6.25
$ws1.Range("S44").Value = 0
$ws1.Range("T44").Value = 39
$ws1.Range("U44").Value = 60.94
$ws1.Range("S45").Value = 0
$ws1.Range("T45").Value = 39
$ws1.Range("U45").Value = 60.94
$ws1.Range("S46").Value = 1
$ws1.Range("T46").Value = 0
$ws1.Range("U46").Value = 0
$ws1.Range("S47").Value = 0
$ws1.Range("T47").Value = 32
$ws1.Range("U47").Value = 50
$ws1.Range("S48").Value = 0
$ws1.Range("T48").Value = 4
$ws1.Range("U48").Value = 6.25
$ws1.Range("S49").Value = 0
$ws1.Range("T49").Value = 36
$ws1.Range("U49").Value = 56.25
$ws1.Range("S50").Value = 0
$ws1.Range("T50").Value = 4
$ws1.Range("U50").Value = 6.25
$ws1.Range("S51").Value = 0
$ws1.Range("T51").Value = 41
$ws1.Range("U51").Value = 64.06
$ws1.Range("S52").Value = 1
$ws1.Range("T52").Value = 0
$ws1.Range("U52").Value = 0
$ws1.Range("S53").Value = 0
$ws1.Range("T53").Value = 39
$ws1.Range("U53").Value = 60.94
$ws1.Range("S54").Value = 0
$ws1.Range("T54").Value = 15
$ws1.Range("U54").Value = 23.44
$ws1.Range("S55").Value = 0
$ws1.Range("T55").Value = 35
$ws1.Range("U55").Value = 54.69
$ws1.Range("S56").Value = 1
$ws1.Range("T56").Value = 0
$ws1.Range("U56").Value = 0
$ws1.Range("S57").Value = 0
$ws1.Range("T57").Value = 8
$ws1.Range("U57").Value = 12.5
$ws1.Range("S58").Value = 0
$ws1.Range("T58").Value = 2
$ws1.Range("U58").Value = 3.13
$ws1.Range("S59").Value = 0
$ws1.Range("T59").Value = 39
$ws1.Range("U59").Value = 60.94
$ws1.Range("S60").Value = 0
$ws1.Range("T60").Value = 39
$ws1.Range("U60").Value = 60.94
$ws1.Range("S61").Value = 1
$ws1.Range("T61").Value = 0
$ws1.Range("U61").Value = 0
$ws1.Range("S62").Value = 0
$ws1.Range("T62").Value = 39
$ws1.Range("U62").Value = 60.94
$ws1.Range("S63").Value = 0
$ws1.Range("T63").Value = 34
$ws1.Range("U63").Value = 53.13
$ws1.Range("S64").Value = 1
$ws1.Range("T64").Value = 0
$ws1.Range("U64").Value = 0
$ws1.Range("S65").Value = 0
$ws1.Range("T65").Value = 39
$ws1.Range("U65").Value = 60.94
$ws1.Range("S66").Value = 0
$ws1.Range("T66").Value = 4
$ws1.Range("U66").Value = 6.25
$ws1.Range("S67").Value = 0
$ws1.Range("T67").Value = 39
$ws1.Range("U67").Value = 60.94
$ws1.Range("S68").Value = 0
$ws1.Range("T68").Value = 18
$ws1.Range("U68").Value = 28.13
$ws1.Range("S69").Value = 1
$ws1.Range("T69").Value = 0
$ws1.Range("U69").Value = 0
$ws1.Range("S70").Value = 0
$ws1.Range("T70").Value = 39
$ws1.Range("U70").Value = 60.94
$ws1.Range("S71").Value = 0
$ws1.Range("T71").Value = 39
$ws1.Range("U71").Value = 60.94
$ws1.Range("S72").Value = 0
$ws1.Range("T72").Value = 39
$ws1.Range("U72").Value = 60.94
$ws1.Range("S73").Value = 1
$ws1.Range("T73").Value = 0
$ws1.Range("U73").Value = 0
$ws1.Range("S74").Value = 0
$ws1.Range("T74").Value = 4
$ws1.Range("U74").Value = 6.25
$ws1.Range("S75").Value = 0
$ws1.Range("T75").Value = 2
$ws1.Range("U75").Value = 3.13
$ws1.Range("S76").Value = 0
$ws1.Range("T76").Value = 9
$ws1.Range("U76").Value = 14.06
$ws1.Range("S77").Value = 0
$ws1.Range("T77").Value = 14
$ws1.Range("U77").Value = 21.88
$ws1.Range("S78").Value = 0
$ws1.Range("T78").Value = 41
$ws1.Range("U78").Value = 64.06
$ws1.Range("S79").Value = 1
$ws1.Range("T79").Value = 0
$ws1.Range("U79").Value = 0
$ws1.Range("S80").Value = 0
$ws1.Range("T80").Value = 39
$ws1.Range("U80").Value = 60.94
$ws1.Range("S81").Value = 0
$ws1.Range("T81").Value = 6
$ws1.Range("U81").Value = 9.38
$ws1.Range("S82").Value = 1
$ws1.Range("T82").Value = 0
$ws1.Range("U82").Value = 0
$ws1.Range("S83").Value = 0
$ws1.Range("T83").Value = 21
$ws1.Range("U83").Value = 32.81
$ws1.Range("S84").Value = 1
$ws1.Range("T84").Value = 0
$ws1.Range("U84").Value = 0
$ws1.Range("S85").Value = 0
$ws1.Range("T85").Value = 28
$ws1.Range("U85").Value = 43.75
$ws1.Range("S86").Value = 1
$ws1.Range("T86").Value = 0
$ws1.Range("U86").Value = 0
$ws1.Range("S87").Value = 1
$ws1.Range("T87").Value = 0
$ws1.Range("U87").Value = 0
$ws1.Range("S88").Value = 1
$ws1.Range("T88").Value = 0
$ws1.Range("U88").Value = 0
$ws1.Range("S89").Value = 0
$ws1.Range("T89").Value = 5
$ws1.Range("U89").Value = 7.81
$ws1.Range("S90").Value = 1
$ws1.Range("T90").Value = 0
$ws1.Range("U90").Value = 0
$ws1.Range("S91").Value = 0
$ws1.Range("T91").Value = 39
$ws1.Range("U91").Value = 60.94
$ws1.Range("S92").Value = 0
$ws1.Range("T92").Value = 9
$ws1.Range("U92").Value = 14.06
$ws1.Range("S93").Value = 0
$ws1.Range("T93").Value = 39
$ws1.Range("U93").Value = 60.94
$ws1.Range("S94").Value = 0
$ws1.Range("T94").Value = 37
$ws1.Range("U94").Value = 57.81
$ws1.Range("S95").Value = 1
$ws1.Range("T95").Value = 0
$ws1.Range("U95").Value = 0
$ws1.Range("S96").Value = 0
$ws1.Range("T96").Value = 27
$ws1.Range("U96").Value = 42.19
$ws1.Range("S97").Value = 1
$ws1.Range("T97").Value = 0
$ws1.Range("U97").Value = 0
$ws1.Range("S98").Value = 1
$ws1.Range("T98").Value = 0
$ws1.Range("U98").Value = 0
$ws1.Range("S99").Value = 0
$ws1.Range("T99").Value = 28
$ws1.Range("U99").Value = 43.75
$ws1.Range("S100").Value = 0
$ws1.Range("T100").Value = 39
$ws1.Range("U100").Value = 60.94
$ws1.Range("S101").Value = 0
$ws1.Range("T101").Value = 2
$ws1.Range("U101").Value = 3.13
$ws1.Range("S102").Value = 0
$ws1.Range("T102").Value = 32
$ws1.Range("U102").Value = 50
$ws1.Range("S103").Value = 0
$ws1.Range("T103").Value = 39
$ws1.Range("U103").Value = 60.94
$ws1.Range("S104").Value = 1
$ws1.Range("T104").Value = 0
$ws1.Range("U104").Value = 0
$ws1.Range("S105").Value = 0
$ws1.Range("T105").Value = 39
$ws1.Range("U105").Value = 60.94
$ws1.Range("S106").Value = 0
$ws1.Range("T106").Value = 2
$ws1.Range("U106").Value = 3.13
$ws1.Range("S107").Value = 1
$ws1.Range("T107").Value = 0
$ws1.Range("U107").Value = 0
$ws1.Range("S108").Value = 1
$ws1.Range("T108").Value = 0
$ws1.Range("U108").Value = 0
$ws1.Range("S109").Value = 0
$ws1.Range("T109").Value = 2
$ws1.Range("U109").Value = 3.13
$ws1.Range("S110").Value = 1
$ws1.Range("T110").Value = 0
$ws1.Range("U110").Value = 0
$ws1.Range("S111").Value = 0
$ws1.Range("T111").Value = 39
$ws1.Range("U111").Value = 60.94
$ws1.Range("S112").Value = 0
$ws1.Range("T112").Value = 34
$ws1.Range("U112").Value = 53.13
$ws1.Range("S113").Value = 0
$ws1.Range("T113").Value = 39
$ws1.Range("U113").Value = 60.94
$ws1.Range("S114").Value = 1
$ws1.Range("T114").Value = 0
$ws1.Range("U114").Value = 0
$ws1.Range("S115").Value = 0
$ws1.Range("T115").Value = 39
$ws1.Range("U115").Value = 60.94
$ws1.Range("S116").Value = 0
$ws1.Range("T116").Value = 2
$ws1.Range("U116").Value = 3.13
$ws1.Range("S117").Value = 0
$ws1.Range("T117").Value = 37
$ws1.Range("U117").Value = 57.81
$ws1.Range("S118").Value = 0
$ws1.Range("T118").Value = 35
$ws1.Range("U118").Value = 54.69
$ws1.Range("S119").Value = 0
$ws1.Range("T119").Value = 7
$ws1.Range("U119").Value = 10.94
$ws1.Range("S120").Value = 0
$ws1.Range("T120").Value = 39
$ws1.Range("U120").Value = 60.94
$ws1.Range("S121").Value = 1
$ws1.Range("T121").Value = 0
$ws1.Range("U121").Value = 0
$ws1.Range("S122").Value = 0
$ws1.Range("T122").Value = 21
$ws1.Range("U122").Value = 32.81
$ws1.Range("S123").Value = 1
$ws1.Range("T123").Value = 0
$ws1.Range("U123").Value = 0
$ws1.Range("S124").Value = 0
$ws1.Range("T124").Value = 6
$ws1.Range("U124").Value = 9.38
$ws1.Range("S125").Value = 1
$ws1.Range("T125").Value = 0
$ws1.Range("U125").Value = 0
$ws1.Range("S126").Value = 0
$ws1.Range("T126").Value = 25
$ws1.Range("U126").Value = 39.06
$ws1.Range("S127").Value = 1
$ws1.Range("T127").Value = 0
$ws1.Range("U127").Value = 0
$ws1.Range("S128").Value = 0
$ws1.Range("T128").Value = 28
$ws1.Range("U128").Value = 43.75
$ws1.Range("S129").Value = 0
$ws1.Range("T129").Value = 39
$ws1.Range("U129").Value = 60.94
$ws1.Range("S130").Value = 0
$ws1.Range("T130").Value = 9
$ws1.Range("U130").Value = 14.06
$ws1.Range("S131").Value = 1
$ws1.Range("T131").Value = 0
$ws1.Range("U131").Value = 0
$ws1.Range("S132").Value = 0
$ws1.Range("T132").Value = 15
$ws1.Range("U132").Value = 23.44
$ws1.Range("S133").Value = 0
$ws1.Range("T133").Value = 37
$ws1.Range("U133").Value = 57.81
$ws1.Range("S134").Value = 1
$ws1.Range("T134").Value = 0
$ws1.Range("U134").Value = 0
$ws1.Range("S135").Value = 0
$ws1.Range("T135").Value = 7
$ws1.Range("U135").Value = 10.94
$ws1.Range("S136").Value = 0
$ws1.Range("T136").Value = 32
$ws1.Range("U136").Value = 50
$ws1.Range("S137").Value = 0
$ws1.Range("T137").Value = 39
$ws1.Range("U137").Value = 60.94
$ws1.Range("S138").Value = 0
$ws1.Range("T138").Value = 17
$ws1.Range("U138").Value = 26.56
$ws1.Range("S139").Value = 0
$ws1.Range("T139").Value = 3
$ws1.Range("U139").Value = 4.69
$ws1.Range("S140").Value = 0
$ws1.Range("T140").Value = 31
$ws1.Range("U140").Value = 48.44
$ws1.Range("S141").Value = 0
$ws1.Range("T141").Value = 8
$ws1.Range("U141").Value = 12.5
$ws1.Range("S142").Value = 1
$ws1.Range("T142").Value = 0
$ws1.Range("U142").Value = 0
$ws1.Range("S143").Value = 0
$ws1.Range("T143").Value = 2
$ws1.Range("U143").Value = 3.13
$ws1.Range("S144").Value = 0
$ws1.Range("T144").Value = 8
$ws1.Range("U144").Value = 12.5
# --- Add new eSGARP AEI column (G) to sheet2 ---
$ws2.Range("G1").Value = "eSGARP"
$ws2.Range("G2").Value = "AEI"

$ws2.Range("G3").Value = 1
$ws2.Range("G4").Value = 0.3333333333332575
$ws2.Range("G5").Value = 0.3333333333332575
$ws2.Range("G6").Value = 0.8333333333334849
$ws2.Range("G7").Value = 0.3333333333332575
$ws2.Range("G8").Value = 1
$ws2.Range("G9").Value = 1
$ws2.Range("G10").Value = 0.3333333333332575
$ws2.Range("G11").Value = 1
$ws2.Range("G12").Value = 1
$ws2.Range("G13").Value = 1
$ws2.Range("G14").Value = 1
$ws2.Range("G15").Value = 1
$ws2.Range("G16").Value = 0.3333333333332575
$ws2.Range("G17").Value = 1
$ws2.Range("G18").Value = 0.8333333333334849
$ws2.Range("G19").Value = 1
$ws2.Range("G20").Value = 1
$ws2.Range("G21").Value = 0.3333333333332575
$ws2.Range("G22").Value = 0.900000000000091
$ws2.Range("G23").Value = 0.3333333333332575
$ws2.Range("G24").Value = 0.3333333333332575
$ws2.Range("G25").Value = 1
$ws2.Range("G26").Value = 0.8333333333334849
$ws2.Range("G27").Value = 1
$ws2.Range("G28").Value = 0.9999999999990905
$ws2.Range("G29").Value = 1
$ws2.Range("G30").Value = 1
$ws2.Range("G31").Value = 0.5000000000004547
$ws2.Range("G32").Value = 0.7500000000004547
$ws2.Range("G33").Value = 0.9999999999990905
$ws2.Range("G34").Value = 0.3333333333332575
$ws2.Range("G35").Value = 0.9999999999990905
$ws2.Range("G36").Value = 0.3333333333332575
$ws2.Range("G37").Value = 0.3333333333332575
$ws2.Range("G38").Value = 1
$ws2.Range("G39").Value = 0.9999999999990905
$ws2.Range("G40").Value = 0.7500000000004547
$ws2.Range("G41").Value = 1
$ws2.Range("G42").Value = 0.7999999999997272
$ws2.Range("G43").Value = 0.9999999999990905
$ws2.Range("G44").Value = 0.3333333333332575
$ws2.Range("G45").Value = 0.3333333333332575
$ws2.Range("G46").Value = 1
$ws2.Range("G47").Value = 0.7000000000002728
$ws2.Range("G48").Value = 0.9999999999990905
$ws2.Range("G49").Value = 0.5000000000004547
$ws2.Range("G50").Value = 0.9999999999990905
$ws2.Range("G51").Value = 0.3333333333332575
$ws2.Range("G52").Value = 1
$ws2.Range("G53").Value = 0.3333333333332575
$ws2.Range("G54").Value = 0.8333333333334849
$ws2.Range("G55").Value = 0.6666666666665151
$ws2.Range("G56").Value = 1
$ws2.Range("G57").Value = 0.9499999999998181
$ws2.Range("G58").Value = 0.9999999999990905
$ws2.Range("G59").Value = 0.3333333333332575
$ws2.Range("G60").Value = 0.3333333333332575
$ws2.Range("G61").Value = 1
$ws2.Range("G62").Value = 0.3333333333332575
$ws2.Range("G63").Value = 0.3333333333332575
$ws2.Range("G64").Value = 1
$ws2.Range("G65").Value = 0.3333333333332575
$ws2.Range("G66").Value = 0.9999999999990905
$ws2.Range("G67").Value = 0.3333333333332575
$ws2.Range("G68").Value = 0.8333333333334849
$ws2.Range("G69").Value = 1
$ws2.Range("G70").Value = 0.3333333333332575
$ws2.Range("G71").Value = 0.3333333333332575
$ws2.Range("G72").Value = 0.3333333333332575
$ws2.Range("G73").Value = 1
$ws2.Range("G74").Value = 0.9749999999994543
$ws2.Range("G75").Value = 0.9999999999990905
$ws2.Range("G76").Value = 0.8750000000004547
$ws2.Range("G77").Value = 0.7999999999997272
$ws2.Range("G78").Value = 0.3999999999998636
$ws2.Range("G79").Value = 1
$ws2.Range("G80").Value = 0.5000000000004547
$ws2.Range("G81").Value = 0.8750000000004547
$ws2.Range("G82").Value = 1
$ws2.Range("G83").Value = 0.7500000000004547
$ws2.Range("G84").Value = 1
$ws2.Range("G85").Value = 0.7500000000004547
$ws2.Range("G86").Value = 1
$ws2.Range("G87").Value = 1
$ws2.Range("G88").Value = 1
$ws2.Range("G89").Value = 0.900000000000091
$ws2.Range("G90").Value = 1
$ws2.Range("G91").Value = 0.3333333333332575
$ws2.Range("G92").Value = 0.933333333333394
$ws2.Range("G93").Value = 0.3333333333332575
$ws2.Range("G94").Value = 0.3333333333332575
$ws2.Range("G95").Value = 1
$ws2.Range("G96").Value = 0.3333333333332575
$ws2.Range("G97").Value = 1
$ws2.Range("G98").Value = 1
$ws2.Range("G99").Value = 0.900000000000091
$ws2.Range("G100").Value = 0.3333333333332575
$ws2.Range("G101").Value = 0.9999999999990905
$ws2.Range("G102").Value = 0.3333333333332575
$ws2.Range("G103").Value = 0.3333333333332575
$ws2.Range("G104").Value = 1
$ws2.Range("G105").Value = 0.3333333333332575
$ws2.Range("G106").Value = 0.9999999999990905
$ws2.Range("G107").Value = 1
$ws2.Range("G108").Value = 1
$ws2.Range("G109").Value = 0.9999999999990905
$ws2.Range("G110").Value = 1
$ws2.Range("G111").Value = 0.3333333333332575
$ws2.Range("G112").Value = 0.599999999999909
$ws2.Range("G113").Value = 0.3333333333332575
$ws2.Range("G114").Value = 1
$ws2.Range("G115").Value = 0.3333333333332575
$ws2.Range("G116").Value = 0.9999999999990905
$ws2.Range("G117").Value = 0.599999999999909
$ws2.Range("G118").Value = 0.3999999999998636
$ws2.Range("G119").Value = 0.900000000000091
$ws2.Range("G120").Value = 0.3333333333332575
$ws2.Range("G121").Value = 1
$ws2.Range("G122").Value = 0.5000000000004547
$ws2.Range("G123").Value = 1
$ws2.Range("G124").Value = 0.9999999999990905
$ws2.Range("G125").Value = 1
$ws2.Range("G126").Value = 0.8333333333334849
$ws2.Range("G127").Value = 1
$ws2.Range("G128").Value = 0.7500000000004547
$ws2.Range("G129").Value = 0.3333333333332575
$ws2.Range("G130").Value = 0.8750000000004547
$ws2.Range("G131").Value = 1
$ws2.Range("G132").Value = 0.8750000000004547
$ws2.Range("G133").Value = 0.599999999999909
$ws2.Range("G134").Value = 1
$ws2.Range("G135").Value = 0.900000000000091
$ws2.Range("G136").Value = 0.6666666666665151
$ws2.Range("G137").Value = 0.3333333333332575
$ws2.Range("G138").Value = 0.8333333333334849
$ws2.Range("G139").Value = 0.8750000000004547
$ws2.Range("G140").Value = 0.3999999999998636
$ws2.Range("G141").Value = 0.3333333333332575
$ws2.Range("G142").Value = 1
$ws2.Range("G143").Value = 0.9999999999990905
$ws2.Range("G144").Value = 0.8333333333334849
